$wb = $excel.ActiveWorkbook

# --- Link sheet: p_nom 100 -> 10, p_nom_extendable False -> True (rows 2 & 3) ---
$wsLink = $wb.Worksheets.Item("Link")

$wsLink.Range("J2").Value = 10
$wsLink.Range("J3").Value = 10

# Copy a cell already containing the literal text "True" within the same sheet
# so the new value keeps its original shared-string text type (not a native boolean).
$wsLink.Range("G2").Copy()
$wsLink.Range("L2").PasteSpecial(-4163)
$wsLink.Range("G3").Copy()
$wsLink.Range("L3").PasteSpecial(-4163)

# --- Store sheet: e_nom 0 -> 100, e_nom_extendable True -> False (row 2) ---
$wsStore = $wb.Worksheets.Item("Store")

$wsStore.Range("E2").Value = 100

$wsStore.Range("M2").Copy()
$wsStore.Range("G2").PasteSpecial(-4163)

# --- Selection / active sheet bookkeeping ---
[void]$wsStore.Range("G3").Select()

$wsLink.Activate()
[void]$wsLink.Range("J4").Select()
